$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value  = 4445.37816546568
$ws.Range("G3").Value  = 3522.860453977319
$ws.Range("G4").Value  = 6675.338893189715
$ws.Range("G5").Value  = 3586.112347731234
$ws.Range("G6").Value  = 0.756964996274725
$ws.Range("G7").Value  = 0.8419859648424165
$ws.Range("G8").Value  = 0.8323486429190367
$ws.Range("G9").Value  = 0.6867279060865942
$ws.Range("G10").Value = 0.3639975625079706
$ws.Range("G11").Value = 0.3944289289741915
$ws.Range("G12").Value = 0.3200883900502076
$ws.Range("G13").Value = 0.4438841847577117
